$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (Exposure period and Notes changed)
$ws.Range("C2").Value = "27/12/20 7.30pm - 9pm"
$ws.Range("D2").Value = "Case dined for dinner"

# Add new row 3
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "Left Bank Melbourne, 1 Southbank Blvd"
$ws.Range("C3").Value = "25/12/20 12pm - 3pm"
$ws.Range("D3").Value = "Case ate in store"
$ws.Range("E3").Value = "new"

# Add new row 4
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C4").Value = "28/12/2020 10pm - 12.30am"
$ws.Range("D4").Value = "Case attended Venue"
$ws.Range("E4").Value = "new"

# Resize columns to fit the new content (mirrors Excel's auto column-fit
# behaviour after the new rows/text were entered)
$ws.Columns.AutoFit()

# Update selection to C2 as in the diff
$ws.Range("C2").Select()
